$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, [string]$val)
    $escaped = $val -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '51.426.33'
Set-TextValue $ws.Range('E2') '  +2.76%  '
Set-TextValue $ws.Range('D3') '2.753.03'
Set-TextValue $ws.Range('E3') '  +2.86%  '
Set-TextValue $ws.Range('E4') '  +0.00%  '
Set-TextValue $ws.Range('D5') '115.67'
Set-TextValue $ws.Range('E5') '  +1.34%  '
Set-TextValue $ws.Range('D6') '331.36'
Set-TextValue $ws.Range('E6') '  +1.38%  '
Set-TextValue $ws.Range('D7') '0.532'
Set-TextValue $ws.Range('E7') '  +0.48%  '
Set-TextValue $ws.Range('D8') '1.00'
Set-TextValue $ws.Range('E8') '  +0.00%  '
Set-TextValue $ws.Range('E9') '  +2.45%  '
Set-TextValue $ws.Range('D10') '41.62'
Set-TextValue $ws.Range('E10') '  +1.24%  '
Set-TextValue $ws.Range('D11') '20.31'
Set-TextValue $ws.Range('E11') '  +0.91%  '
Set-TextValue $ws.Range('D12') '0.0830'
Set-TextValue $ws.Range('E12') '  +0.44%  '
Set-TextValue $ws.Range('E13') '  +2.76%  '
Set-TextValue $ws.Range('D14') '7.70'
Set-TextValue $ws.Range('E14') '  +4.14%  '
Set-TextValue $ws.Range('D15') '3.175.79'
Set-TextValue $ws.Range('E15') '  +2.72%  '
Set-TextValue $ws.Range('D16') '2.736.80'
Set-TextValue $ws.Range('E16') '  +1.35%  '
Set-TextValue $ws.Range('D17') '0.887'
Set-TextValue $ws.Range('E17') '  +1.01%  '
Set-TextValue $ws.Range('D18') '51.347.07'
Set-TextValue $ws.Range('E18') '  +2.71%  '
Set-TextValue $ws.Range('D19') '13.62'
Set-TextValue $ws.Range('E19') '  +2.78%  '
Set-TextValue $ws.Range('D20') '3.03'
Set-TextValue $ws.Range('E20') '  +4.36%  '
Set-TextValue $ws.Range('D21') '6.86'
Set-TextValue $ws.Range('E21') '  +0.89%  '
Set-TextValue $ws.Range('D22') '0.0₃0965'
Set-TextValue $ws.Range('E22') '  +0.24%  '
Set-TextValue $ws.Range('D23') '281.95'
Set-TextValue $ws.Range('E23') '  +1.34%  '
Set-TextValue $ws.Range('D24') '70.25'
Set-TextValue $ws.Range('E24') '  -3.31%  '
Set-TextValue $ws.Range('D26') '26.89'
Set-TextValue $ws.Range('E26') '  -0.22%  '
Set-TextValue $ws.Range('D27') '1.00'
Set-TextValue $ws.Range('E27') '  +0.04%  '
Set-TextValue $ws.Range('D28') '10.34'
Set-TextValue $ws.Range('E28') '  +2.11%  '
Set-TextValue $ws.Range('E29') '  -0.37%  '
Set-TextValue $ws.Range('E30') '  -1.47%  '
Set-TextValue $ws.Range('D31') '35.67'
Set-TextValue $ws.Range('E31') '  -2.00%  '
Set-TextValue $ws.Range('D32') '50.17'
Set-TextValue $ws.Range('E32') '  -0.31%  '
Set-TextValue $ws.Range('D33') '5.66'
Set-TextValue $ws.Range('E33') '  +2.85%  '
Set-TextValue $ws.Range('D34') '0.0824'
Set-TextValue $ws.Range('E34') '  +0.55%  '
Set-TextValue $ws.Range('D35') '19.47'
Set-TextValue $ws.Range('E35') '  -1.11%  '
Set-TextValue $ws.Range('E36') '  -0.28%  '
Set-TextValue $ws.Range('D37') '2.11'
Set-TextValue $ws.Range('E37') '  +1.15%  '
Set-TextValue $ws.Range('D38') '5.03'
Set-TextValue $ws.Range('E38') '  -1.82%  '
Set-TextValue $ws.Range('D39') '3.23'
Set-TextValue $ws.Range('E39') '  +1.69%  '
Set-TextValue $ws.Range('D40') '129.27'
Set-TextValue $ws.Range('E40') '  +3.46%  '
Set-TextValue $ws.Range('D41') '23.71'
Set-TextValue $ws.Range('E41') '  +4.16%  '
Set-TextValue $ws.Range('E42') '  +10.39%  '
Set-TextValue $ws.Range('D43') '2.31'
Set-TextValue $ws.Range('E43') '  +4.05%  '
Set-TextValue $ws.Range('E44') '  +0.13%  '
Set-TextValue $ws.Range('D45') '3.43'
Set-TextValue $ws.Range('E45') '  +3.47%  '
Set-TextValue $ws.Range('D46') '2.113.37'
Set-TextValue $ws.Range('E46') '  -0.29%  '
Set-TextValue $ws.Range('D47') '2.24'
Set-TextValue $ws.Range('E47') '  +10.13%  '
Set-TextValue $ws.Range('D48') '2.25'
Set-TextValue $ws.Range('E48') '  -0.47%  '
Set-TextValue $ws.Range('D49') '5.56'
Set-TextValue $ws.Range('E49') '  +2.92%  '
Set-TextValue $ws.Range('D50') '9.08'
Set-TextValue $ws.Range('E50') '  +0.03%  '
Set-TextValue $ws.Range('B51') 'MinaProtocolToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
Set-TextValue $ws.Range('D51') '1.52'
Set-TextValue $ws.Range('E51') '  +7.85%  '

$excel.CutCopyMode = 0
